# Adds the season record columns (Wins, Losses, Ties) to the worksheet,
# matching the style already used by the other header cells in row 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1) onto the three
# new header cells so they pick up the same (bold / bordered / centered) style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-39: season record values (Wins=80, Losses=82, Ties=0)
for ($row = 2; $row -le 39; $row++) {
    $ws.Cells.Item($row, 30).Value = 80
    $ws.Cells.Item($row, 31).Value = 82
    $ws.Cells.Item($row, 32).Value = 0
}
